# Processing corrected to exclude rogue subject names.
# The sheet contains a "-- as percentages" summary row at the end of each
# gender block (rows 68 and 135) that isn't an actual subject and must be
# removed from the data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 135 first (the later row) so that the row index of the
# earlier row (68) is unaffected by the deletion.
$ws.Rows.Item(135).EntireRow.Delete()
$ws.Rows.Item(68).EntireRow.Delete()
